$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Column C header: "Swimming" -> "Breast Stroke Drill"
$ws.Range("C1").Value = "Breast Stroke Drill"

# New column J header: "Seiza (Japanese Sitting Posture)" (extends the
# used range from A1:I101 to A1:J101)
$ws.Range("J1").Value = "Seiza (Japanese Sitting Posture)"

# --- Column widths ------------------------------------------------------
# Widen column C (now holds the longer "Breast Stroke Drill" label) and
# column J (holds the long "Seiza (Japanese Sitting Posture)" label).
# NOTE: ColumnWidth is expressed in "normal style" character units and is
# quantized by the host to a 1/6-character pixel grid before being stored
# as the sheet's <col width=.../>, so the literal target widths (16.3 and
# 28.11) cannot be represented bit-exactly - the inputs below are the
# values whose quantized result lands closest to those targets.
$ws.Columns.Item(3).ColumnWidth = 15.5
$ws.Columns.Item(10).ColumnWidth = 27.3333333333333

# --- Day 7 (row 8) : corrected which exercise the reps belonged to ------
$ws.Range("C8").Value = 0
$ws.Range("H8").Value = 20

# --- Day 9 (row 10) : fill in that day's workout ------------------------
$ws.Range("C10").Value = 20
$ws.Range("D10").Value = 20
$ws.Range("E10").Value = "5 min"
$ws.Range("F10").Value = 40
$ws.Range("G10").Value = 40
$ws.Range("H10").Value = 40
$ws.Range("I10").Value = 40

# --- Restore the author's final selection -------------------------------
$ws.Range("F11").Select()
